$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New B-column values for rows 2..23 (A holds 0..21 respectively)
$values = @(
    45.27692569068709,
    53.15072906367325,
    43.01162633521314,
    50.99019513592785,
    25.49509756796392,
    35,
    46.09772228646444,
    49.49747468305833,
    39.05124837953327,
    40.31128874149275,
    47.43416490252569,
    30.4138126514911,
    36.40054944640259,
    55.90169943749475,
    52.20153254455276,
    49.24428900898052,
    55.22680508593631,
    31.6227766016838,
    25,
    54.08326913195985,
    38.07886552931955,
    53.85164807134504
)

$lastExistingRow = 16
$firstNewRow = $lastExistingRow + 1
$lastNewRow = 2 + $values.Length - 1

# Extend column A's formatting (style + border + alignment) to the newly
# added rows by copying the format of the last existing styled cell.
$ws.Range("A" + $lastExistingRow).Copy()
$ws.Range("A" + $firstNewRow + ":A" + $lastNewRow).PasteSpecial(-4122)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $i
    $ws.Cells.Item($row, 2).Value = $values[$i]
}
